$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------------
# Row 424: day-separator row ("THURSDAY"), formatted like the existing
# separator rows (e.g. row 5).
# ---------------------------------------------------------------------------
$ws.Range("A5:F5").Copy()
$ws.Range("A424:F424").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(424, 2).Value = "THURSDAY"

# ---------------------------------------------------------------------------
# Rows 425-428: plain data rows (no explicit row height), formatted like the
# regular data rows (e.g. row 2).
# ---------------------------------------------------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A425:F428").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(425, 1).Value = "AV Shutdown"
$ws.Cells.Item(425, 2).Value = 42656
$ws.Cells.Item(425, 3).Value = "1900"
$ws.Cells.Item(425, 4).Value = "LSB"
$ws.Cells.Item(425, 5).Value = "103"
$ws.Cells.Item(425, 6).Value = "Make sure neck mic goes back to drawer and log off touchscreen."

$ws.Cells.Item(426, 1).Value = "AV Shutdown"
$ws.Cells.Item(426, 2).Value = 42656
$ws.Cells.Item(426, 3).Value = "1900"
$ws.Cells.Item(426, 4).Value = "LSB"
$ws.Cells.Item(426, 5).Value = "105"
$ws.Cells.Item(426, 6).Value = "Make sure neck mic goes back to drawer and log off touchscreen."

$ws.Cells.Item(427, 1).Value = "AV Shutdown"
$ws.Cells.Item(427, 2).Value = 42656
$ws.Cells.Item(427, 3).Value = "1730"
$ws.Cells.Item(427, 4).Value = "LSB"
$ws.Cells.Item(427, 5).Value = "107"
$ws.Cells.Item(427, 6).Value = "Make sure neck mic goes back to drawer and log off touchscreen."

$ws.Cells.Item(428, 1).Value = "AV Shutdown"
$ws.Cells.Item(428, 2).Value = 42656
$ws.Cells.Item(428, 3).Value = "1730"
$ws.Cells.Item(428, 4).Value = "LSB"
$ws.Cells.Item(428, 5).Value = "101"
$ws.Cells.Item(428, 6).Value = "Make sure neck mic goes back to drawer and log off touchscreen."

# ---------------------------------------------------------------------------
# Row 429: taller data row (Special Instructions wraps to 3 lines), formatted
# like the existing ht=45 data rows (e.g. row 6).
# ---------------------------------------------------------------------------
$ws.Range("A6:F6").Copy()
$ws.Range("A429:F429").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(429).RowHeight = 45

$ws.Cells.Item(429, 1).Value = "AV Shutdown"
$ws.Cells.Item(429, 2).Value = 42656
$ws.Cells.Item(429, 3).Value = "1900"
$ws.Cells.Item(429, 4).Value = "CLH"
$ws.Cells.Item(429, 5).Value = "L"
$ws.Cells.Item(429, 6).Value = "PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN'T GET LOGGED OFF."

# ---------------------------------------------------------------------------
# Row 430: taller data row (Special Instructions wraps to 2 lines), formatted
# like the existing ht=30 data rows (e.g. row 20).
# ---------------------------------------------------------------------------
$ws.Range("A20:F20").Copy()
$ws.Range("A430:F430").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(430).RowHeight = 30

$ws.Cells.Item(430, 1).Value = "Lockup"
$ws.Cells.Item(430, 2).Value = 42656
$ws.Cells.Item(430, 3).Value = "2030"
$ws.Cells.Item(430, 4).Value = "CLH"
$ws.Cells.Item(430, 5).Value = "K"
$ws.Cells.Item(430, 6).Value = "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS."

# ---------------------------------------------------------------------------
# Reflect the new selection/active cell left behind by the edit.
# ---------------------------------------------------------------------------
[void]$ws.Range("F434").Select()

Write-Host "Appended rows 424-430 to the Logs sheet."
